# Week 17 data logged for the Falcons Players Data workbook.
# Sheet "Rushing" (index 1) and "Receiving" (index 2) both get updated
# counting stats, and a new player (F.Darby) is added to the Receiving
# sheet, pushing the remaining receiving rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Rushing
# ---------------------------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

# M.Ryan (row 2)
$rushing.Cells.Item(2,4).Value = 11   # D2 2DATT
$rushing.Cells.Item(2,6).Value = 7    # F2 RZATT

# M.Davis (row 4)
$rushing.Cells.Item(4,3).Value = 90   # C4 1DATT
$rushing.Cells.Item(4,4).Value = 34   # D4 2DATT
$rushing.Cells.Item(4,6).Value = 18   # F4 RZATT

# C.Patterson (row 5)
$rushing.Cells.Item(5,3).Value = 101  # C5 1DATT
$rushing.Cells.Item(5,4).Value = 38   # D5 2DATT
$rushing.Cells.Item(5,6).Value = 33   # F5 RZATT

# Q.Ollison (row 7)
$rushing.Cells.Item(7,3).Value = 15   # C7 1DATT
$rushing.Cells.Item(7,4).Value = 3    # D7 2DATT

# ---------------------------------------------------------------------
# Sheet 2: Receiving
# ---------------------------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

# M.Davis (row 2)
$receiving.Cells.Item(2,3).Value = 53  # C2 Short Target
$receiving.Cells.Item(2,4).Value = 39  # D2 Short Comp

# C.Patterson (row 3)
$receiving.Cells.Item(3,3).Value = 59  # C3 Short Target
$receiving.Cells.Item(3,4).Value = 44  # D3 Short Comp

# R.Gage (row 6)
$receiving.Cells.Item(6,3).Value = 76  # C6 Short Target
$receiving.Cells.Item(6,4).Value = 64  # D6 Short Comp
$receiving.Cells.Item(6,5).Value = 16  # E6 Deep Target
$receiving.Cells.Item(6,6).Value = 10  # F6 Deep Comp
$receiving.Cells.Item(6,7).Value = 11  # G6 RZ Target
$receiving.Cells.Item(6,8).Value = 7   # H6 RZ Comp

# O.Zaccheaus (row 7)
$receiving.Cells.Item(7,3).Value = 42  # C7 Short Target
$receiving.Cells.Item(7,4).Value = 24  # D7 Short Comp

# C.Blake (row 8)
$receiving.Cells.Item(8,5).Value = 2   # E8 Deep Target
$receiving.Cells.Item(8,7).Value = 2   # G8 RZ Target

# Push rows 10-15 down to 11-16 to make room for the new player F.Darby,
# copying the formatted ID cell (column A) so the new rows inherit the
# same style as the rest of the table.
$receiving.Cells.Item(15,1).Copy($receiving.Cells.Item(16,1))
$receiving.Cells.Item(14,1).Copy($receiving.Cells.Item(15,1))
$receiving.Cells.Item(13,1).Copy($receiving.Cells.Item(14,1))
$receiving.Cells.Item(12,1).Copy($receiving.Cells.Item(13,1))
$receiving.Cells.Item(11,1).Copy($receiving.Cells.Item(12,1))
$receiving.Cells.Item(10,1).Copy($receiving.Cells.Item(11,1))

# Row 16: K.Smith (was row 15)
$receiving.Cells.Item(16,1).Value = 14
$receiving.Cells.Item(16,2).Value = "K.Smith"
$receiving.Cells.Item(16,3).Value = 5
$receiving.Cells.Item(16,4).Value = 5
$receiving.Cells.Item(16,5).Value = 0
$receiving.Cells.Item(16,6).Value = 0
$receiving.Cells.Item(16,7).Value = 0
$receiving.Cells.Item(16,8).Value = 0

# Row 15: P.Hesse (was row 14)
$receiving.Cells.Item(15,2).Value = "P.Hesse"
$receiving.Cells.Item(15,3).Value = 5
$receiving.Cells.Item(15,4).Value = 4
$receiving.Cells.Item(15,5).Value = 0
$receiving.Cells.Item(15,6).Value = 0
$receiving.Cells.Item(15,7).Value = 0
$receiving.Cells.Item(15,8).Value = 0

# Row 14: L.Smith (was row 13)
$receiving.Cells.Item(14,2).Value = "L.Smith"
$receiving.Cells.Item(14,3).Value = 10
$receiving.Cells.Item(14,4).Value = 8
$receiving.Cells.Item(14,5).Value = 1
$receiving.Cells.Item(14,6).Value = 1
$receiving.Cells.Item(14,7).Value = 2
$receiving.Cells.Item(14,8).Value = 1

# Row 13: H.Hurst (was row 12)
$receiving.Cells.Item(13,2).Value = "H.Hurst"
$receiving.Cells.Item(13,3).Value = 29
$receiving.Cells.Item(13,4).Value = 25
$receiving.Cells.Item(13,5).Value = 1
$receiving.Cells.Item(13,6).Value = 0
$receiving.Cells.Item(13,7).Value = 7
$receiving.Cells.Item(13,8).Value = 6

# Row 12: K.Pitts (was row 11)
$receiving.Cells.Item(12,2).Value = "K.Pitts"
$receiving.Cells.Item(12,3).Value = 76
$receiving.Cells.Item(12,4).Value = 50
$receiving.Cells.Item(12,5).Value = 28
$receiving.Cells.Item(12,6).Value = 16
$receiving.Cells.Item(12,7).Value = 14
$receiving.Cells.Item(12,8).Value = 5

# Row 11: M.Hall (was row 10)
$receiving.Cells.Item(11,2).Value = "M.Hall"
$receiving.Cells.Item(11,3).Value = 1
$receiving.Cells.Item(11,4).Value = 1
$receiving.Cells.Item(11,5).Value = 0
$receiving.Cells.Item(11,6).Value = 0
$receiving.Cells.Item(11,7).Value = 0
$receiving.Cells.Item(11,8).Value = 0

# Row 10: F.Darby (brand-new player row)
$receiving.Cells.Item(10,2).Value = "F.Darby"
$receiving.Cells.Item(10,3).Value = 0
$receiving.Cells.Item(10,4).Value = 0
$receiving.Cells.Item(10,5).Value = 1
$receiving.Cells.Item(10,6).Value = 0
$receiving.Cells.Item(10,7).Value = 0
$receiving.Cells.Item(10,8).Value = 0

Write-Output "Week 17 data logged."
